$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from 2023-10-25 (45224) to 2023-11-03 (45233)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 45233
}
